$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.482.84'
$ws.Range("E2").Value = '  +9.34%  '

# Row 3
$ws.Range("D3").Value = '1.611.54'
$ws.Range("E3").Value = '  +9.34%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.67%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.23'
$ws.Range("E5").Value = '  +9.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9901'
$ws.Range("E6").Value = '  +4.23%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3683'
$ws.Range("E7").Value = '  +1.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3401'
$ws.Range("E8").Value = '  +11.39%  '

# Row 9
$ws.Range("E9").Value = '  +7.11%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.137'
$ws.Range("E10").Value = '  +7.89%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07065'
$ws.Range("E11").Value = '  +6.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.31%  '

# Row 13
$ws.Range("E13").Value = '  +9.93%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.915'
$ws.Range("E14").Value = '  +7.38%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.650'
$ws.Range("E15").Value = '  +7.24%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.607.45'
$ws.Range("E16").Value = '  +9.19%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001085'
$ws.Range("E17").Value = '  +5.56%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9911'
$ws.Range("E18").Value = '  +4.29%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06694'
$ws.Range("E19").Value = '  +12.68%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.21'
$ws.Range("E20").Value = '  +12.91%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.018'
$ws.Range("E21").Value = '  +9.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.06'
$ws.Range("E22").Value = '  +11.32%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.84'
$ws.Range("E23").Value = '  +6.94%  '

# Row 24
$ws.Range("D24").Value = '22.540.97'
$ws.Range("E24").Value = '  +9.52%  '

# Row 25
$ws.Range("E25").Value = '  +4.95%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.561'
$ws.Range("E26").Value = '  +21.50%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.59'
$ws.Range("E27").Value = '  +4.67%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.58'
$ws.Range("E28").Value = '  +13.87%  '

# Row 29
$ws.Range("D29").Value = '1.790.63'
$ws.Range("E29").Value = '  +9.67%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.53'
$ws.Range("E30").Value = '  +7.91%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.059'
$ws.Range("E31").Value = '  +2.74%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.165'
$ws.Range("E32").Value = '  +23.26%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9542'
$ws.Range("E33").Value = '  +18.39%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.683'
$ws.Range("E34").Value = '  +11.66%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08244'
$ws.Range("E35").Value = '  +3.71%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.97'
$ws.Range("E36").Value = '  +15.68%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.250'
$ws.Range("E37").Value = '  +11.43%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.270'
$ws.Range("E38").Value = '  +4.62%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.574'
$ws.Range("E39").Value = '  +16.01%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06100'
$ws.Range("E40").Value = '  +4.30%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02216'
$ws.Range("E41").Value = '  +8.29%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2022'
$ws.Range("E42").Value = '  +7.85%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9899'
$ws.Range("E43").Value = '  +4.03%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5920'
$ws.Range("E44").Value = '  +11.92%  '

# Row 45
$ws.Range("E45").Value = '  +8.28%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.23'
$ws.Range("E46").Value = '  +8.54%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5699'
$ws.Range("E47").Value = '  +10.04%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.65'
$ws.Range("E48").Value = '  +8.44%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.970'
$ws.Range("E49").Value = '  +8.67%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06835'
$ws.Range("E50").Value = '  +5.65%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.78'
$ws.Range("E51").Value = '  +9.63%  '
